# HW3_Eng.pptx edit: on slide 5 ("Example"), the black caption box that
# shows the example answer count ("文字方塊 10") is updated from "22" to
# "26". The shape uses spAutoFit, so PowerPoint automatically grows its
# height to fit the new text once the run is edited.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)

# Locate the "文字方塊 10" shape (currently holds the text "22") instead of
# hard-coding an index, in case shape ordering ever differs.
$target = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.TextRange.Text -eq "22") {
        $target = $shp
        break
    }
}

if ($target -eq $null) {
    $target = $s.Shapes.Item("文字方塊 10")
}

$tr = $target.TextFrame.TextRange

# Replace just the second character ("2" -> "6") so the run that already
# reads "2" is left untouched and only the trailing digit becomes a new
# run, matching how the text was actually edited by hand.
$tr.Characters(2, 1).Text = "6"
